$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    heading paragraph (paragraph 1).
# ------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

$boldLabel = "Meta description"
$restOfLine = ": Experience the rich Chinese tradition with impeccable graphics and multiple bonuses of 8 Dragons. Play free and choose your free spin options."
$fullMeta = $boldLabel + $restOfLine

$metaRange = $metaPara.Range
$metaRange.Text = $fullMeta

$metaStart = $d.Paragraphs(2).Range.Start
$boldRange = $d.Range($metaStart, $metaStart + $boldLabel.Length)
$boldRange.Bold = 1

# ------------------------------------------------------------------
# 2) Remove the now-duplicated bold "Play 8 Dragons Free..." paragraph
#    that used to sit right before the closing italic paragraph.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupHeading = $d.Paragraphs($count - 1)
$dupHeading.Range.Delete()

# ------------------------------------------------------------------
# 3) Replace the text of the trailing italic paragraph with the new
#    image-prompt copy (keep its italic formatting intact).
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$lastStart = $lastPara.Range.Start
$lastEnd = $lastPara.Range.End
$lastRange = $d.Range($lastStart, $lastEnd)

$newPrompt = 'Create a feature image for "8 Dragons" that features a happy Maya warrior with glasses. The image should be in a cartoon style and should have a vibrant and eye-catching color scheme. The Maya warrior should be depicted holding a dragon in one hand and a pile of gold coins in the other, surrounded by Chinese-themed symbols such as lanterns and scrolls. In the background, you can add a colorful dragon or a temple to add to the overall theme of the game. The image should convey the excitement and adventure of playing "8 Dragons" and entice players to try their luck at this exciting slot game.'

$lastRange.Text = $newPrompt
